# Auto-generated edit script: updates Leve profit/price figures across multiple sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW) per the scheduled market-price refresh.

$wb = $excel.ActiveWorkbook

$sheetEdits = @{}

$sheetEdits["ALC"] = @{
    "H2" = 579.4
    "I2" = 99.25
    "K2" = 99.25
    "M2" = 13.75
    "H74" = 11873.75
    "I74" = 10997.667
    "J74" = 14502
    "K74" = 10997.667
    "L74" = 14502
    "M74" = -10061.667
    "N74" = -16374
    "H77" = 11873.75
    "I77" = 10997.667
    "J77" = 14502
    "K77" = 54988.335
    "L77" = 72510
    "M77" = -50308.335
    "N77" = -81870
    "H113" = 12097.857
    "J113" = 15727
    "L113" = 15727
    "N113" = -22235
    "H132" = 1512.88
    "I132" = 1512.88
    "J132" = 0
    "K132" = 4538.64
    "L132" = 0
    "M132" = -2008.64
    "N132" = "DELETE"
}

$sheetEdits["ARM"] = @{
    "H61" = 3827.0435
    "I61" = 3101.0454
    "J61" = 19799
    "K61" = 3101.0454
    "L61" = 19799
    "M61" = -2889.0454
    "N61" = -20223
    "H62" = 15000
    "J62" = 15000
    "L62" = 15000
    "N62" = -16248
    "H65" = 15000
    "J65" = 15000
    "L65" = 45000
    "N65" = -51240
    "H110" = 8526.352999999999
    "I110" = 7729
    "K110" = 7729
    "M110" = -5684
    "H132" = 1926.6888
    "I132" = 1107.6052
    "J132" = 6373.143
    "K132" = 3322.8156
    "L132" = 19119.429
    "M132" = -792.8155999999999
    "N132" = -24179.429
    "H136" = 3827.0435
    "I136" = 3101.0454
    "J136" = 19799
    "K136" = 9303.136200000001
    "L136" = 59397
    "M136" = -6753.136200000001
    "N136" = -64497
}

$sheetEdits["BSM"] = @{
    "H5" = 1574.625
    "I5" = 118.333336
    "J5" = 2448.4
    "K5" = 118.333336
    "L5" = 2448.4
    "M5" = -5.333336000000003
    "N5" = -2674.4
    "H10" = 2083.1667
    "I10" = 1919.8
    "J10" = 2900
    "K10" = 1919.8
    "L10" = 2900
    "M10" = -1779.8
    "N10" = -3180
    "H86" = 10914.857
    "I86" = 6424.8
    "K86" = 6424.8
    "M86" = -5301.8
    "H89" = 10914.857
    "I89" = 6424.8
    "K89" = 32124
    "M89" = -26508
    "H107" = 596.36365
    "I107" = 566
    "K107" = 566
    "M107" = 1354
    "H118" = 0
    "J118" = 0
    "L118" = 0
    "N118" = "DELETE"
    "H134" = 3185.742
    "I134" = 1619.28
    "J134" = 9712.666999999999
    "K134" = 4857.84
    "L134" = 29138.001
    "M134" = -2322.84
    "N134" = -34208.001
}

$sheetEdits["CRP"] = @{
    "H86" = 14428.143
    "I86" = 9497
    "J86" = 16400.6
    "K86" = 9497
    "L86" = 16400.6
    "M86" = -8374
    "N86" = -18646.6
    "H89" = 14428.143
    "I89" = 9497
    "J89" = 16400.6
    "K89" = 47485
    "L89" = 82003
    "M89" = -41869
    "N89" = -93235
    "H99" = 2202.8
    "I99" = 2202.8
    "K99" = 2202.8
    "M99" = -704.8000000000002
    "H107" = 991.46155
    "I107" = 895.4545000000001
    "J107" = 1519.5
    "K107" = 895.4545000000001
    "L107" = 1519.5
    "M107" = 1024.5455
    "N107" = -5359.5
    "H126" = 2202.8
    "I126" = 2202.8
    "K126" = 6608.400000000001
    "M126" = -4138.400000000001
    "H132" = 3751.9375
    "I132" = 3414.5862
    "J132" = 7013
    "K132" = 10243.7586
    "L132" = 21039
    "M132" = -7713.758600000001
    "N132" = -26099
    "H134" = 2381.5938
    "I134" = 1713.7241
    "J134" = 8837.666999999999
    "K134" = 5141.1723
    "L134" = 26513.001
    "M134" = -2606.1723
    "N134" = -31583.001
}

$sheetEdits["CUL"] = @{
    "H107" = 756.6667
    "I107" = 756.6667
    "K107" = 2270.0001
    "M107" = -350.0001000000002
    "H140" = 2903.2222
    "I140" = 2423.0625
    "J140" = 6744.5
    "K140" = 7269.1875
    "L140" = 20233.5
    "M140" = -2089.1875
    "N140" = -30593.5
}

$sheetEdits["GSM"] = @{
    "H32" = 44944.5
    "J32" = 44944.5
    "L32" = 44944.5
    "N32" = -45536.5
    "H80" = 9633.817999999999
    "I80" = 9069
    "J80" = 10622.25
    "K80" = 9069
    "L80" = 10622.25
    "M80" = -8071
    "N80" = -12618.25
    "H83" = 9633.817999999999
    "I83" = 9069
    "J83" = 10622.25
    "K83" = 45345
    "L83" = 53111.25
    "M83" = -40353
    "N83" = -63095.25
    "H113" = 11656.333
    "I113" = 9984.5
    "J113" = 15000
    "K113" = 9984.5
    "L113" = 15000
    "M113" = -7814.5
    "N113" = -19340
    "H132" = 26744.234
    "I132" = 36695.227
    "J132" = 7464.1875
    "K132" = 110085.681
    "L132" = 22392.5625
    "M132" = -107555.681
    "N132" = -27452.5625
}

$sheetEdits["LTW"] = @{
    "H43" = 20057.295
    "I43" = 16749.75
    "K43" = 16749.75
    "M43" = -16556.75
    "H62" = 295000
    "J62" = 295000
    "L62" = 295000
    "N62" = -296248
    "H65" = 295000
    "J65" = 295000
    "L65" = 885000
    "N65" = -891240
    "H82" = 3682.1304
    "J82" = 6551.273
    "L82" = 6551.273
    "N82" = -7273.273
    "H85" = 3682.1304
    "J85" = 6551.273
    "L85" = 6551.273
    "N85" = -9047.273000000001
    "H93" = 1613.2106
    "I93" = 1185.125
    "K93" = 1185.125
    "M93" = 62.875
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetEdits[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $newValue = $cellMap[$cellRef]
        if ("$newValue" -eq "DELETE") {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $newValue
        }
    }
}

Write-Host "Applied $($sheetEdits.Values | ForEach-Object { $_.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell edits across $($sheetEdits.Keys.Count) sheets."